# Updated cryptos list with latest price/volume data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.845.61"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.76%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.140.40"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.74%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "571.41"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.88%  "
$ws.Range("E6").Value = "  +3.61%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.137.35"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.75%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.526"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.80%  "
$ws.Range("E10").Value = "  +4.78%  "
$ws.Range("E11").Value = "  +0.95%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.501"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +6.22%  "
$ws.Range("E13").Value = "  +10.79%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "37.47"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +6.39%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.652.50"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.94%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.910.53"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.83%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.20"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +6.23%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.138.65"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.68%  "
$ws.Range("E19").Value = "  +0.24%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "511.29"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +6.48%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.92"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +7.16%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.733"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +8.71%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "15.36"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +11.35%  "
$ws.Range("E24").Value = "  +3.20%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "84.99"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.61%  "
$ws.Range("E26").Value = "  +0.41%  "
$ws.Range("E27").Value = "  +3.34%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.70"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +8.35%  "
$ws.Range("E29").Value = "  +4.88%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "27.92"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +6.37%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.999"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.01%  "
$ws.Range("E32").Value = "  +3.65%  "
$ws.Range("E33").Value = "  +5.58%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.04"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +7.78%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.57"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +5.75%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "55.40"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.73%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "477.97"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +4.23%  "
$ws.Range("E38").Value = "  +3.74%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0860"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.26%  "
$ws.Range("E40").Value = "  -1.33%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.115.34"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.59%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.61"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.13%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.120"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.83%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.292"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +11.48%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.47"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +14.86%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "29.09"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.70%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0₃0572"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +10.44%  "
$ws.Range("E49").Value = "  +3.07%  "
$ws.Range("E50").Value = "  +10.33%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "118.67"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.95%  "
